# Update the multiplication problems in the table to the new values
# described in the commit "Update master to output generated at c8c62b6".
$d = $word.ActiveDocument

$replacements = @(
    @("745×5=", "969×5="),
    @("837×9=", "830×5="),
    @("151×5=", "296×8="),
    @("406×2=", "988×8="),
    @("206×4=", "586×5="),
    @("680×6=", "534×7="),
    @("668×3=", "362×2="),
    @("179×3=", "353×7="),
    @("734×5=", "902×6="),
    @("495×5=", "810×3="),
    @("155×4=", "911×9="),
    @("286×6=", "543×5="),
    @("661×6=", "768×6="),
    @("832×2=", "543×8="),
    @("354×9=", "726×4="),
    @("963×3=", "621×5="),
    @("905×9=", "960×6="),
    @("367×8=", "687×7="),
    @("464×8=", "882×9="),
    @("631×9=", "563×2="),
    @("449×5=", "568×9="),
    @("820×8=", "588×7="),
    @("345×2=", "865×2="),
    @("255×7=", "704×4="),
    @("869×6=", "287×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
